# Apply "fixed workflow" results update to both sheets (NBR, BAR).
# For rows 2-16 (Cutoff = A = 0..14), update the Reaction (B) and
# Reaction_number (C) columns with the new values, then remove the
# trailing 4 rows (17-20) that no longer exist in the new results.

$wb = $excel.ActiveWorkbook

# New B (Reaction) / C (Reaction_number) values for Cutoff 0..14, per sheet.
$nbrValues = @(
    @(5, 839),
    @(6, 0),
    @(7, 846),
    @(8, 835),
    @(9, 824),
    @(10, 819),
    @(11, 819),
    @(12, 812),
    @(13, 807),
    @(14, 806),
    @(15, 809),
    @(16, 792),
    @(17, 778),
    @(18, 778),
    @(19, 777)
)

$barValues = @(
    @(5, 751),
    @(6, 0),
    @(7, 731),
    @(8, 731),
    @(9, 725),
    @(10, 722),
    @(11, 721),
    @(12, 714),
    @(13, 710),
    @(14, 709),
    @(15, 703),
    @(16, 694),
    @(17, 692),
    @(18, 692),
    @(19, 692)
)

$sheetData = @{
    "NBR" = $nbrValues
    "BAR" = $barValues
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    $values = $sheetData[$name]

    if ($values) {
        for ($i = 0; $i -lt $values.Count; $i++) {
            $row = $i + 2
            $pair = $values[$i]
            $ws.Cells.Item($row, 2).Value = $pair[0]
            $ws.Cells.Item($row, 3).Value = $pair[1]
        }

        # The refreshed results only contain 15 data rows (Cutoff 0..14),
        # so drop the now-stale rows 17-20 (previously Cutoff 15..18).
        $ws.Range("A17:A20").EntireRow.Delete() | Out-Null
    }
}

$wb.Save() | Out-Null
